$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 1528.7
$ws.Range("I19").Value = 4440.6
$ws.Range("J19").Value = 558.06665
$ws.Range("K19").Value = 4440.6
$ws.Range("L19").Value = 558.06665
$ws.Range("M19").Value = -4265.6
$ws.Range("N19").Value = -908.06665
# row 129
$ws.Range("H129").Value = 734.61536
$ws.Range("J129").Value = 900
$ws.Range("L129").Value = 2700
$ws.Range("N129").Value = -12700
# row 137
$ws.Range("H137").Value = 27029404
$ws.Range("I137").Value = 1216.3914
$ws.Range("J137").Value = 71432856
$ws.Range("K137").Value = 3649.1742
$ws.Range("L137").Value = 214298568
$ws.Range("M137").Value = -1099.1742
$ws.Range("N137").Value = -214303668
# row 141
$ws.Range("H141").Value = 1188.4
$ws.Range("I141").Value = 1235.5
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 3706.5
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = 1473.5
$ws.Range("N141").Value = -13360

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 2676396.8
$ws.Range("I2").Value = 2399.8333
$ws.Range("J2").Value = 5885193
$ws.Range("K2").Value = 2399.8333
$ws.Range("L2").Value = 5885193
$ws.Range("M2").Value = -2286.8333
$ws.Range("N2").Value = -5885419
# row 63
$ws.Range("H63").Value = 2834.25
$ws.Range("I63").Value = 2062.3333
$ws.Range("J63").Value = 3465.818
$ws.Range("K63").Value = 2062.3333
$ws.Range("L63").Value = 3465.818
$ws.Range("M63").Value = -1376.3333
$ws.Range("N63").Value = -4837.818
# row 66
$ws.Range("H66").Value = 2834.25
$ws.Range("I66").Value = 2062.3333
$ws.Range("J66").Value = 3465.818
$ws.Range("K66").Value = 10311.6665
$ws.Range("L66").Value = 17329.09
$ws.Range("M66").Value = -6879.666499999999
$ws.Range("N66").Value = -24193.09
# row 116
$ws.Range("H116").Value = 2676396.8
$ws.Range("I116").Value = 2399.8333
$ws.Range("J116").Value = 5885193
$ws.Range("K116").Value = 2399.8333
$ws.Range("L116").Value = 5885193
$ws.Range("M116").Value = -105.8332999999998
$ws.Range("N116").Value = -5889781
# row 122
$ws.Range("H122").Value = 1643.5
$ws.Range("I122").Value = 1391.3334
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 4174.0002
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -1724.0002
$ws.Range("N122").Value = -12100
# row 132
$ws.Range("H132").Value = 1530996.1
$ws.Range("I132").Value = 2779070.8
$ws.Range("K132").Value = 8337212.399999999
$ws.Range("M132").Value = -8334682.399999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 2676396.8
$ws.Range("I3").Value = 2399.8333
$ws.Range("J3").Value = 5885193
$ws.Range("K3").Value = 2399.8333
$ws.Range("L3").Value = 5885193
$ws.Range("M3").Value = -2285.8333
$ws.Range("N3").Value = -5885421
# row 86
$ws.Range("H86").Value = 1300.5938
$ws.Range("I86").Value = 1260.76
$ws.Range("J86").Value = 1442.8572
$ws.Range("K86").Value = 1260.76
$ws.Range("L86").Value = 1442.8572
$ws.Range("M86").Value = -137.76
$ws.Range("N86").Value = -3688.8572
# row 89
$ws.Range("H89").Value = 1300.5938
$ws.Range("I89").Value = 1260.76
$ws.Range("J89").Value = 1442.8572
$ws.Range("K89").Value = 6303.8
$ws.Range("L89").Value = 7214.286
$ws.Range("M89").Value = -687.8000000000002
$ws.Range("N89").Value = -18446.286

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3530.6333
$ws.Range("I31").Value = 1972.5
$ws.Range("J31").Value = 4097.227
$ws.Range("K31").Value = 1972.5
$ws.Range("L31").Value = 4097.227
$ws.Range("M31").Value = -1677.5
$ws.Range("N31").Value = -4687.227
# row 34
$ws.Range("H34").Value = 3530.6333
$ws.Range("I34").Value = 1972.5
$ws.Range("J34").Value = 4097.227
$ws.Range("K34").Value = 1972.5
$ws.Range("L34").Value = 4097.227
$ws.Range("M34").Value = -1770.5
$ws.Range("N34").Value = -4501.227
# row 132
$ws.Range("H132").Value = 1734.9138
$ws.Range("I132").Value = 1620.7709
$ws.Range("J132").Value = 2282.8
$ws.Range("K132").Value = 4862.3127
$ws.Range("L132").Value = 6848.400000000001
$ws.Range("M132").Value = -2332.3127
$ws.Range("N132").Value = -11908.4
# row 134
$ws.Range("H134").Value = 1954.9149
$ws.Range("I134").Value = 1803.5143
$ws.Range("J134").Value = 2396.5
$ws.Range("K134").Value = 5410.5429
$ws.Range("L134").Value = 7189.5
$ws.Range("M134").Value = -2875.5429
$ws.Range("N134").Value = -12259.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 100
$ws.Range("H100").Value = 2607.3333
$ws.Range("J100").Value = 2709.75
$ws.Range("L100").Value = 8129.25
$ws.Range("N100").Value = -9751.25
# row 117
$ws.Range("H117").Value = 727.6429000000001
$ws.Range("I117").Value = 448.5
$ws.Range("J117").Value = 1099.8334
$ws.Range("K117").Value = 1345.5
$ws.Range("L117").Value = 3299.5002
$ws.Range("M117").Value = 2096.5
$ws.Range("N117").Value = -10183.5002
# row 121
$ws.Range("H121").Value = 22727862
$ws.Range("I121").Value = 581.0714
$ws.Range("J121").Value = 62500604
$ws.Range("K121").Value = 1743.2142
$ws.Range("L121").Value = 187501812
$ws.Range("M121").Value = -433.2142000000001
$ws.Range("N121").Value = -187504432
# row 125
$ws.Range("H125").Value = 2098.5715
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 2281.6667
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 6845.000100000001
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -16685.0001
# row 129
$ws.Range("H129").Value = 2830.0645
$ws.Range("I129").Value = 685.5714
$ws.Range("J129").Value = 3455.5417
$ws.Range("K129").Value = 2056.7142
$ws.Range("L129").Value = 10366.6251
$ws.Range("M129").Value = 2943.2858
$ws.Range("N129").Value = -20366.6251
# row 131
$ws.Range("H131").Value = 4123.6665
$ws.Range("J131").Value = 2453.258
$ws.Range("L131").Value = 7359.773999999999
$ws.Range("N131").Value = -17439.774

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 49
$ws.Range("H49").Value = 7000
$ws.Range("J49").Value = 7000
$ws.Range("L49").Value = 7000
$ws.Range("N49").Value = -7368

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 1519.909
$ws.Range("I40").Value = 1471.9
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1471.9
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1335.9
$ws.Range("N40").Value = -2272
# row 42
$ws.Range("H42").Value = 5545.4546
# row 49
$ws.Range("H49").Value = 5545.4546

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 1181136.6
$ws.Range("I132").Value = 1595384
$ws.Range("J132").Value = 2124.8462
$ws.Range("K132").Value = 4786152
$ws.Range("L132").Value = 6374.5386
$ws.Range("M132").Value = -4783622
$ws.Range("N132").Value = -11434.5386
